{"js": "// Highlight quantitative metrics (percentages, dollar amounts, large numbers)\n// in specific resume bullet paragraphs with bold + color (#2C3E50), matching\n// the commit \"Implement quantitative metrics highlighting across all resume\n// formats\".\n//\n// Strategy: locate each target paragraph by its exact original text, then\n// within that paragraph search for each metric token (in left-to-right\n// order) and apply bold + color formatting to just that sub-run. Scoping the\n// search to the specific paragraph (rather than the whole document) avoids\n// accidentally reformatting identical-looking numbers that appear elsewhere\n// in the resume (e.g. \"23% to 64%\" also appears in the summary and in the\n// Key Projects section, but only the Professional Experience bullet should\n// change).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Each entry: the paragraph's exact original text, and the ordered list of\n// metric substrings within it that should become bold + colored.\nconst edits = [\n  {\n    text: \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    tokens: [\"23%\", \"64%\"]\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00b14.2% to \\u00b12.1%\",\n    tokens: [\"87%\", \"71%\", \"\\u00b14.2%\", \"\\u00b12.1%\"]\n  },\n  {\n    text: \"\\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    tokens: [\"1,200\"]\n  },\n  {\n    text: \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    tokens: [\"$400M\", \"$1B\"]\n  },\n  {\n    text: \"\\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    tokens: [\"73.5%\", \"$4.7M\"]\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    tokens: [\"87%\", \"71%\"]\n  }\n];\n\n// Track how many paragraphs we've already matched for each distinct text, so\n// that duplicate paragraph texts (none currently, but defensive) each match\n// a distinct paragraph instance rather than reusing the first one found.\nconst usedIndices = new Set();\n\nfor (const edit of edits) {\n  let target = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (usedIndices.has(i)) continue;\n    if (paragraphs.items[i].text === edit.text) {\n      target = paragraphs.items[i];\n      usedIndices.add(i);\n      break;\n    }\n  }\n  if (!target) {\n    throw new Error(\"Could not locate target paragraph: \" + edit.text);\n  }\n\n  for (const token of edit.tokens) {\n    const hits = target.search(token, { matchCase: true });\n    hits.load(\"items\");\n    await context.sync();\n    if (hits.items.length === 0) {\n      throw new Error(\"Could not locate token '\" + token + \"' in paragraph: \" + edit.text);\n    }\n    const hit = hits.items[0];\n    hit.font.bold = true;\n    hit.font.color = \"#2C3E50\";\n  }\n  await context.sync();\n}\n", "ps1": "# Highlight quantitative metrics (percentages, dollar amounts, large numbers)\n# in specific resume bullet paragraphs with bold + color (#2C3E50), matching\n# the commit \"Implement quantitative metrics highlighting across all resume\n# formats\".\n#\n# Strategy: walk every paragraph in the document, match it against the exact\n# original bullet text, and for each metric token inside that paragraph use a\n# paragraph-scoped Find.Execute to locate just that sub-string, then flip\n# Font.Bold/Font.Color on the found (collapsed) range. Scoping Find to the\n# paragraph's own Range (rather than $d.Content) avoids reformatting\n# identical-looking numbers elsewhere in the resume (e.g. \"23% to 64%\" also\n# appears in the summary and in the Key Projects section, but only the\n# Professional Experience bullet should change).\n\n$d = $word.ActiveDocument\n\n$edits = @(\n    @{\n        Text   = [char]8226 + \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Tokens = @(\"23%\", \"64%\")\n    },\n    @{\n        Text   = [char]8226 + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \" + [char]177 + \"4.2% to \" + [char]177 + \"2.1%\"\n        Tokens = @(\"87%\", \"71%\", ([char]177 + \"4.2%\"), ([char]177 + \"2.1%\"))\n    },\n    @{\n        Text   = [char]8226 + \" Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        Tokens = @(\"1,200\")\n    },\n    @{\n        Text   = [char]8226 + \" Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        Tokens = @(\"`$400M\", \"`$1B\")\n    },\n    @{\n        Text   = [char]8226 + \" Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        Tokens = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        Text   = [char]8226 + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        Tokens = @(\"87%\", \"71%\")\n    }\n)\n\n$usedParas = @{}\n$count = $d.Paragraphs.Count\n\nforeach ($edit in $edits) {\n    $target = $null\n    for ($i = 1; $i -le $count; $i++) {\n        if ($usedParas.ContainsKey($i)) { continue }\n        $p = $d.Paragraphs.Item($i)\n        # Range.Text carries the trailing paragraph mark, so compare against\n        # text + \"`r\".\n        if ($p.Range.Text -eq ($edit.Text + \"`r\")) {\n            $target = $p\n            $usedParas[$i] = $true\n            break\n        }\n    }\n    if ($null -eq $target) {\n        throw \"Could not locate target paragraph: $($edit.Text)\"\n    }\n\n    foreach ($token in $edit.Tokens) {\n        $rng = $target.Range\n        $found = $rng.Find.Execute($token, $true)\n        if (-not $found) {\n            throw \"Could not locate token '$token' in paragraph: $($edit.Text)\"\n        }\n        $rng.Font.Bold = $true\n        $rng.Font.Color = \"#2C3E50\"\n    }\n}\n"}
